# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2210   (columns A:J)
#   *_new -> *_FV2304   (columns L:U, "diff" in K stays as-is)
# Then wrap the data range in a table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
$leftCols  = "A", "B", "C", "D", "E", "F", "G", "H", "I", "J"
$rightCols = "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U"

for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Range($leftCols[$i] + "1").Value = $fields[$i] + "_FV2210"
    $ws.Range($rightCols[$i] + "1").Value = $fields[$i] + "_FV2304"
}

# Turn the used range into an Excel table (ListObject), headers already in row 1.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U59"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
